$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2" '34.521.75'
Set-TextValue "E2" '  +0.25%  '
Set-TextValue "D3" '1.814.66'
Set-TextValue "E3" '  +1.16%  '
Set-TextValue "E4" '  -0.05%  '
Set-TextValue "D5" '229.16'
Set-TextValue "E5" '  +0.72%  '
Set-TextValue "D6" '0.581'
Set-TextValue "E6" '  +4.64%  '
Set-TextValue "E7" '  -0.11%  '
Set-TextValue "D8" '36.24'
Set-TextValue "E8" '  +10.09%  '
Set-TextValue "D9" '0.305'
Set-TextValue "E9" '  +3.19%  '
Set-TextValue "D10" '0.0701'
Set-TextValue "E10" '  +1.28%  '
Set-TextValue "D11" '0.0964'
Set-TextValue "E11" '  +1.87%  '
Set-TextValue "D12" '2.076.37'
Set-TextValue "E12" '  +1.06%  '
Set-TextValue "D13" '11.61'
Set-TextValue "E13" '  +4.24%  '
Set-TextValue "D14" '1.812.67'
Set-TextValue "E14" '  +1.45%  '
Set-TextValue "D15" '0.653'
Set-TextValue "E15" '  +2.43%  '
Set-TextValue "D16" '4.57'
Set-TextValue "E16" '  +6.39%  '
Set-TextValue "D17" '34.521.90'
Set-TextValue "E17" '  +0.08%  '
Set-TextValue "D18" '69.51'
Set-TextValue "E18" '  +0.97%  '
Set-TextValue "D19" '248.90'
Set-TextValue "E19" '  +1.21%  '
Set-TextValue "D20" '0.0₃0803'
Set-TextValue "E20" '  +0.10%  '
Set-TextValue "D21" '11.64'
Set-TextValue "E21" '  +2.94%  '
Set-TextValue "E22" '  +0.08%  '
Set-TextValue "D23" '4.25'
Set-TextValue "E23" '  +2.08%  '
Set-TextValue "D24" '172.62'
Set-TextValue "E24" '  +1.71%  '
Set-TextValue "E25" '  +3.79%  '
Set-TextValue "D26" '8.10'
Set-TextValue "E26" '  +9.68%  '
Set-TextValue "D27" '17.01'
Set-TextValue "E27" '  +2.43%  '
Set-TextValue "E28" '  +4.06%  '
Set-TextValue "E29" '  -0.09%  '
Set-TextValue "D30" '4.11'
Set-TextValue "E30" '  +2.38%  '
Set-TextValue "D31" '3.92'
Set-TextValue "E31" '  +2.98%  '
Set-TextValue "D32" '0.0536'
Set-TextValue "E32" '  +1.82%  '
Set-TextValue "E33" '  +1.12%  '
Set-TextValue "D34" '1.87'
Set-TextValue "E34" '  +1.95%  '
Set-TextValue "D35" '1.405.96'
Set-TextValue "E35" '  -1.13%  '
Set-TextValue "D36" '0.686'
Set-TextValue "E36" '  -0.02%  '
Set-TextValue "D37" '2.54'
Set-TextValue "E37" '  -0.65%  '
Set-TextValue "D38" '1.08'
Set-TextValue "E38" '  -0.02%  '
Set-TextValue "E39" '  +0.82%  '
Set-TextValue "D40" '84.56'
Set-TextValue "E40" '  -0.01%  '
Set-TextValue "D41" '0.975'
Set-TextValue "E41" '  +3.04%  '
Set-TextValue "D42" '2.83'
Set-TextValue "E42" '  +1.69%  '
Set-TextValue "E43" '  +0.32%  '
Set-TextValue "D44" '1.20'
Set-TextValue "E44" '  +8.49%  '
Set-TextValue "D45" '13.54'
Set-TextValue "E45" '  -3.18%  '
Set-TextValue "D46" '6.07'
Set-TextValue "E46" '  -1.24%  '
Set-TextValue "D47" '0.0508'
Set-TextValue "E47" '  -3.39%  '
Set-TextValue "D48" '1.975.53'
Set-TextValue "E48" '  +1.00%  '
Set-TextValue "D49" '106.19'
Set-TextValue "E49" '  +0.94%  '
Set-TextValue "E50" '  +0.01%  '
Set-TextValue "D51" '0.0₆0129'
Set-TextValue "E51" '  -0.21%  '
